$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed values in existing rows 2-7 ---
$ws.Range("B2").Value = "NSE:INDIACEM"
$ws.Range("C2").Value = "NSE:AGARIND"
$ws.Range("D2").Value = "NSE:GODREJPROP"
$ws.Range("E2").Value = "NSE:BAJAJ-AUTO"
$ws.Range("C3").Value = "NSE:AXISBANK"
$ws.Range("E3").Value = "NSE:IDFCFIRSTB"
$ws.Range("C4").Value = "NSE:AXSENSEX"
$ws.Range("E4").Value = "NSE:IIFL"
$ws.Range("C5").Value = "NSE:BALAMINES"
$ws.Range("C6").Value = "NSE:BSLNIFTY"
$ws.Range("C7").Value = "NSE:DEN"

# --- Clear cells that became empty ---
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()

# --- Add new rows 8-31: column A index (sequential) and column C ticker ---
# Copy formatting of A2 (bordered, bold, centered) down through new A cells
$ws.Range("A2").Copy()
$ws.Range("A8:A31").PasteSpecial(-4122)

$ws.Range("A8").Value = 6
$ws.Range("C8").Value = "NSE:DHANUKA"
$ws.Range("A9").Value = 7
$ws.Range("C9").Value = "NSE:EXIDEIND"
$ws.Range("A10").Value = 8
$ws.Range("C10").Value = "NSE:GAIL"
$ws.Range("A11").Value = 9
$ws.Range("C11").Value = "NSE:GLAND"
$ws.Range("A12").Value = 10
$ws.Range("C12").Value = "NSE:HAPPSTMNDS"
$ws.Range("A13").Value = 11
$ws.Range("C13").Value = "NSE:HDFCGROWTH"
$ws.Range("A14").Value = 12
$ws.Range("C14").Value = "NSE:HDFCLOWVOL"
$ws.Range("A15").Value = 13
$ws.Range("C15").Value = "NSE:HDFCNIFTY"
$ws.Range("A16").Value = 14
$ws.Range("C16").Value = "NSE:HDFCSENSEX"
$ws.Range("A17").Value = 15
$ws.Range("C17").Value = "NSE:HEADSUP"
$ws.Range("A18").Value = 16
$ws.Range("C18").Value = "NSE:IDEA"
$ws.Range("A19").Value = 17
$ws.Range("C19").Value = "NSE:IIFL"
$ws.Range("A20").Value = 18
$ws.Range("C20").Value = "NSE:IPCALAB"
$ws.Range("A21").Value = 19
$ws.Range("C21").Value = "NSE:JINDWORLD"
$ws.Range("A22").Value = 20
$ws.Range("C22").Value = "NSE:JTLIND"
$ws.Range("A23").Value = 21
$ws.Range("C23").Value = "NSE:KRISHANA"
$ws.Range("A24").Value = 22
$ws.Range("C24").Value = "NSE:LGBBROSLTD"
$ws.Range("A25").Value = 23
$ws.Range("C25").Value = "NSE:LLOYDSME"
$ws.Range("A26").Value = 24
$ws.Range("C26").Value = "NSE:MANAPPURAM"
$ws.Range("A27").Value = 25
$ws.Range("C27").Value = "NSE:NIFTYBEES"
$ws.Range("A28").Value = 26
$ws.Range("C28").Value = "NSE:PATANJALI"
$ws.Range("A29").Value = 27
$ws.Range("C29").Value = "NSE:PITTIENG"
$ws.Range("A30").Value = 28
$ws.Range("C30").Value = "NSE:PNB"
$ws.Range("A31").Value = 29
$ws.Range("C31").Value = "NSE:SAFARI"

Write-Output "done"
